$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05182466666666666
$ws.Range("H2").Value = 0.155474
$ws.Range("M2").Value = 0.74396
$ws.Range("N2").Value = 2.23188
$ws.Range("O2").Value = 0.006259003216804254
$ws.Range("P2").Value = 0.006259003216804255
$ws.Range("Q2").Value = 0.03855547901333333
$ws.Range("R2").Value = 0.34699931112
$ws.Range("S2").Value = 0.006259003216804254
$ws.Range("T2").Value = 0.006259003216804255

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05182466666666666
$ws.Range("H3").Value = 0.155474
$ws.Range("M3").Value = 88.14978533333333
$ws.Range("N3").Value = 264.449356
$ws.Range("O3").Value = 0.7416121699579786
$ws.Range("P3").Value = 0.7416121699579786
$ws.Range("Q3").Value = 4.568333241638221
$ws.Range("R3").Value = 41.11499917474399
$ws.Range("S3").Value = 0.7416121699579786
$ws.Range("T3").Value = 0.7416121699579786

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05182466666666666
$ws.Range("H4").Value = 0.155474
$ws.Range("M4").Value = 29.76859933333333
$ws.Range("N4").Value = 89.305798
$ws.Range("O4").Value = 0.2504459365921425
$ws.Range("P4").Value = 0.2504459365921425
$ws.Range("Q4").Value = 1.542747737583555
$ws.Range("R4").Value = 13.884729638252
$ws.Range("S4").Value = 0.2504459365921425
$ws.Range("T4").Value = 0.2504459365921425

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.05182466666666666
$ws.Range("H5").Value = 0.155474
$ws.Range("M5").Value = 0.2000323333333334
$ws.Range("N5").Value = 0.6000970000000001
$ws.Range("O5").Value = 0.00168289023307462
$ws.Range("P5").Value = 0.00168289023307462
$ws.Range("Q5").Value = 0.01036660899755556
$ws.Range("R5").Value = 0.09329948097800002
$ws.Range("S5").Value = 0.00168289023307462
$ws.Range("T5").Value = 0.00168289023307462
